$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow edits, then restore protection afterward.
$ws.Unprotect()

# Update the confidential disclosure date (2021-05-07 -> 2021-05-10)
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05715775081964215
$ws.Range("E2").Value = -0.01174857142857144

$ws.Range("D3").Value = 0.02384164102188913
$ws.Range("E3").Value = -0.005178365937859697

$ws.Range("D4").Value = 0.03145309373675165
$ws.Range("E4").Value = -0.01869158878504673

$ws.Range("D5").Value = 0.03315084224395451
$ws.Range("E5").Value = -0.002620741295394913

$ws.Range("D6").Value = 0.03957014758457877
$ws.Range("E6").Value = 0.00693343898573695

$ws.Range("D7").Value = 0.01958910436090775
$ws.Range("E7").Value = -0.0060690943043884

$ws.Range("D8").Value = 0.004319166267651131
$ws.Range("E8").Value = -0.02994555353902006

$ws.Range("D9").Value = 0.006966060054179015
$ws.Range("E9").Value = 0.0003750937734434689

$ws.Range("D10").Value = 0.07137206690918974
$ws.Range("E10").Value = 0.002745744096650293

$ws.Range("D11").Value = 0.07141126079546604
$ws.Range("E11").Value = 0.003293084522502987

$ws.Range("D12").Value = 0.1455190609667433
$ws.Range("E12").Value = -0.01012712777418656

$ws.Range("D13").Value = 0.3814211833810083
$ws.Range("E13").Value = -0.002008909075028265

$ws.Range("D14").Value = 0.1142286218580384
$ws.Range("E14").Value = -0.001441095232376455

$ws.Range("E15").Value = -0.003414440726105972

# Restore sheet protection
$ws.Protect("D382", $true, $true, $true, $true)
